$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "Judgment" column (C),
# shifting it to column E, and giving room for the two new columns.
$ws.Range("C:D").EntireColumn.Insert()

# Match the 50-character width used by the other columns.
$ws.Columns.Item(3).ColumnWidth = 49.17
$ws.Columns.Item(4).ColumnWidth = 49.17

# Header row: new column titles.
$ws.Range("C1").Value = "Prosecution Counsel Statement"
$ws.Range("D1").Value = "Defense Counsel Statement"

# Row 2: rewritten/condensed scenario text in column A.
$ws.Range("A2").Value = "** On 26.4.2017 at 1130 hours, police officials were returning to the police station after a search operation when they were attacked by six unknown terrorists near a graveyard. The terrorists opened fire and threw hand grenades, injuring two constables. The police returned fire. Four terrorists fled, while two remained and continued firing near a tube-well. Army personnel, QRF and CTD staff arrived, and the two remaining terrorists blew themselves up. Police recovered two 9mm pistols, live cartridges, and a hand grenade from the terrorists, and took possession of three motorcycles. FIR No.13 was registered under Sections 324/353/148/149 PPC, Sections 3/4/5 of the Explosive Substances Act, Section 7 of the Anti-Terrorism Act, 1997, and Section 15 of the Arms Act at Police Station C.T.D, D.I.Khan. Rahmatullah and Rafiullah were arrested in connection with the case while already in custody for another case (FIR No.461 dated 22.10.2016 under Section 15 Arms Act)."

# Clear the remaining columns on row 2 (Witnesses, new Prosecution/Defense
# Statement columns, and the shifted Judgment column) since the old
# Witnesses/Judgment text no longer applies.
$ws.Range("B2:E2").ClearContents()
